$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "https://www.latimes.com/"
$ws.Range("D2").Style = "Hyperlink"
$ws.Range("A3").Select()
